# Scen_NCAP_NAS.xlsx edit
# Insert a new "LO / CAP_BND / 2045 / 18 / ELE_NEW_WIND-OFF" row into the
# "Nowe moce w śmieciach" block (after the offshore-wind rows), shifting
# every following row down by one, and correct two values in the
# "Nowe moce w jądrówce" (nuclear) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15 (pushes old rows 15..43 down to 16..44),
# inheriting formatting from the row above exactly like Excel's own
# "Insert Sheet Rows" command.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15.
$ws.Range("B15").Value = "LO"
$ws.Range("C15").Value = "CAP_BND"
$ws.Range("D15").Value = 2045
$ws.Range("E15").Value = 18
$ws.Range("F15").Value = "ELE_NEW_WIND-OFF"

# Correct the nuclear (ELE_NEW_NUC_PWR) capacity bound values, now living
# at rows 23 (year 2040) and 24 (year 2050) after the shift.
$ws.Range("E23").Value = 4.4
$ws.Range("E24").Value = 7.7

# Update the sheet's view state to match: scrolled so row 14 is the first
# visible row, with F24 as the active/selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("F24").Select()
